$wb = $excel.ActiveWorkbook

# --- Sheet1 ("crop") ---
$ws1 = $wb.Worksheets.Item("crop")

# Update existing row 2 values
$ws1.Cells.Item(2, 12).Value = 12          # L2: 14 -> 12
$ws1.Cells.Item(2, 18).Value = "unknown"   # R2: 12 -> "unknown" (text)

# Widen column C
$ws1.Columns.Item(3).ColumnWidth = 14.5

# New row 3: Lettuce
$ws1.Cells.Item(3, 1).Value = "BIA project_crop overview 20220316"
$ws1.Cells.Item(3, 3).Value = "Lettuce"
$ws1.Cells.Item(3, 4).Value = "unknown"
$ws1.Cells.Item(3, 5).Value = "unknown"
$ws1.Cells.Item(3, 6).Value = 16
$ws1.Cells.Item(3, 7).Value = 18
$ws1.Cells.Item(3, 8).Value = "unknown"
$ws1.Cells.Item(3, 9).Value = "unknown"
$ws1.Cells.Item(3, 10).Value = 30
$ws1.Cells.Item(3, 11).Value = 45
$ws1.Cells.Item(3, 12).Value = 13
$ws1.Cells.Item(3, 13).Value = 7
$ws1.Cells.Item(3, 14).Value = 7
$ws1.Cells.Item(3, 15).Value = 0.5
$ws1.Cells.Item(3, 16).Value = 0.7
$ws1.Cells.Item(3, 17).Value = 14.4
$ws1.Cells.Item(3, 18).Value = "unknown"
$ws1.Cells.Item(3, 19).Value = 2
$ws1.Cells.Item(3, 20).Value = 4
$ws1.Cells.Item(3, 21).Value = 326
$ws1.Cells.Item(3, 22).Value = 278
$ws1.Cells.Item(3, 23).Value = 240
$ws1.Cells.Item(3, 24).Value = 218
$ws1.Rows.Item(3).RowHeight = 68
$ws1.Cells.Item(3, 1).WrapText = $true

# New row 4: LettuceRomaine
$ws1.Cells.Item(4, 1).Value = "BIA project_crop overview 20220316"
$ws1.Cells.Item(4, 3).Value = "LettuceRomaine"
$ws1.Cells.Item(4, 4).Value = "unknown"
$ws1.Cells.Item(4, 5).Value = "unknown"
$ws1.Cells.Item(4, 6).Value = 16
$ws1.Cells.Item(4, 7).Value = 18
$ws1.Cells.Item(4, 8).Value = "unknown"
$ws1.Cells.Item(4, 9).Value = "unknown"
$ws1.Cells.Item(4, 10).Value = 30
$ws1.Cells.Item(4, 11).Value = 45
$ws1.Cells.Item(4, 12).Value = 13
$ws1.Cells.Item(4, 13).Value = 7
$ws1.Cells.Item(4, 14).Value = 7
$ws1.Cells.Item(4, 15).Value = 0.5
$ws1.Cells.Item(4, 16).Value = 0.7
$ws1.Cells.Item(4, 17).Value = 14.4
$ws1.Cells.Item(4, 18).Value = "unknown"
$ws1.Cells.Item(4, 19).Value = 2
$ws1.Cells.Item(4, 20).Value = 4
$ws1.Cells.Item(4, 21).Value = 241
$ws1.Cells.Item(4, 22).Value = 343
$ws1.Cells.Item(4, 23).Value = 413
$ws1.Cells.Item(4, 24).Value = 330
$ws1.Rows.Item(4).RowHeight = 68
$ws1.Cells.Item(4, 1).WrapText = $true

# --- Sheet2 ("cost") ---
$ws2 = $wb.Worksheets.Item("cost")

$ws2.Columns.Item(1).ColumnWidth = 22.333333333333332

# New row 3: Lettuce
$ws2.Cells.Item(3, 1).Value = "Lettuce"
$ws2.Cells.Item(3, 2).Value = 15
$ws2.Cells.Item(3, 3).Value = 3
$ws2.Cells.Item(3, 4).Value = 10
$ws2.Cells.Item(3, 5).Value = 84.6
$ws2.Cells.Item(3, 6).Value = 42.5
$ws2.Cells.Item(3, 7).Value = 1.24
$ws2.Cells.Item(3, 8).Value = 4.11
$ws2.Cells.Item(3, 9).Value = 0.43

# New row 4: LettuceRomaine
$ws2.Cells.Item(4, 1).Value = "LettuceRomaine"
$ws2.Cells.Item(4, 2).Value = 12.8
$ws2.Cells.Item(4, 3).Value = 3
$ws2.Cells.Item(4, 4).Value = 10
$ws2.Cells.Item(4, 5).Value = 84.6
$ws2.Cells.Item(4, 6).Value = 42.5
$ws2.Cells.Item(4, 7).Value = 1.24
$ws2.Cells.Item(4, 8).Value = 4.11
$ws2.Cells.Item(4, 9).Value = 0.43

# --- Restore selections / active sheet (env must remain the active tab) ---
$ws1.Range("X5").Select()
$ws2.Range("B8").Select()

$ws3 = $wb.Worksheets.Item("env")
$ws3.Range("B1").Select()
